$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.805.82"
$ws.Range("E2").Value = "  -5.62%  "
$ws.Range("D3").Value = "3.273.30"
$ws.Range("E3").Value = "  -6.98%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.53"
$ws.Range("E5").Value = "  -6.52%  "
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.71"
$ws.Range("E6").Value = "  -14.76%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.599"
$ws.Range("E7").Value = "  -1.64%  "
$ws.Range("D8").Value = "3.269.72"
$ws.Range("E8").Value = "  -6.84%  "
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.601"
$ws.Range("E10").Value = "  -8.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.35"
$ws.Range("E11").Value = "  -11.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.131"
$ws.Range("E12").Value = "  -8.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000255"
$ws.Range("E13").Value = "  -6.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.96"
$ws.Range("E14").Value = "  -9.48%  "
$ws.Range("D15").Value = "3.807.28"
$ws.Range("E15").Value = "  -6.51%  "
$ws.Range("D16").Value = "3.282.82"
$ws.Range("E16").Value = "  -6.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.115"
$ws.Range("E17").Value = "  -6.73%  "
$ws.Range("D18").Value = "63.794.15"
$ws.Range("E18").Value = "  -5.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.22"
$ws.Range("E19").Value = "  -7.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.94"
$ws.Range("E20").Value = "  -7.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.943"
$ws.Range("E21").Value = "  -8.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "369.83"
$ws.Range("E22").Value = "  -5.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.73"
$ws.Range("E23").Value = "  -6.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.62"
$ws.Range("E24").Value = "  -4.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.86"
$ws.Range("E25").Value = "  -11.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.84"
$ws.Range("E26").Value = "  -1.63%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.08"
$ws.Range("E27").Value = "  -1.65%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.63"
$ws.Range("E28").Value = "  -7.01%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.23"
$ws.Range("E29").Value = "  -8.30%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.22"
$ws.Range("E30").Value = "  -7.35%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.45"
$ws.Range("E31").Value = "  -8.54%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "636.68"
$ws.Range("E32").Value = "  -11.15%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.56"
$ws.Range("E33").Value = "  -7.55%  "
$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.11"
$ws.Range("E34").Value = "  -5.56%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.104"
$ws.Range("E35").Value = "  -6.81%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.53"
$ws.Range("E36").Value = "  -8.47%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.11"
$ws.Range("E38").Value = "  -6.25%  "
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.382"
$ws.Range("E39").Value = "  -4.10%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0688"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.123"
$ws.Range("E42").Value = "  -6.77%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.896.66"
$ws.Range("E43").Value = "  -5.25%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.42"
$ws.Range("E44").Value = "  -6.62%  "
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.65"
$ws.Range("E45").Value = "  -11.97%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.64"
$ws.Range("E46").Value = "  -4.38%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0391"
$ws.Range("E47").Value = "  -3.57%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.04"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.124"
$ws.Range("E49").Value = "  -2.64%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.73"
$ws.Range("E50").Value = "  +3.85%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.24"
$ws.Range("E51").Value = "  -3.23%  "
